$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M3").Value = 1.02
$ws.Range("N4").Value = 5.9
$ws.Range("L4").Value = 1.01
$ws.Range("Q4").Value = 1.44
$ws.Range("W4").Value = 1.38
$ws.Range("V4").Value = 1.83
$ws.Range("M5").Value = 1.06
$ws.Range("O5").Value = 1.3
$ws.Range("N5").Value = 3.8
$ws.Range("L5").Value = 1.01
$ws.Range("S5").Value = 3.2
$ws.Range("AH5").Value = 18
$ws.Range("AL5").Value = 55
$ws.Range("AK5").Value = 46
$ws.Range("W5").Value = 1.32
$ws.Range("U5").Value = 2.14
$ws.Range("AF5").Value = 29
$ws.Range("Z5").Value = 14
$ws.Range("AG5").Value = 16
$ws.Range("AA5").Value = 26
$ws.Range("AO5").Value = 15.5
$ws.Range("AM5").Value = 110
$ws.Range("AC5").Value = 8.6
$ws.Range("AN5").Value = 46
$ws.Range("T5").Value = 1.75
$ws.Range("AI5").Value = 36
$ws.Range("X5").Value = 19
$ws.Range("R5").Value = 1.38
$ws.Range("Y5").Value = 10.5
$ws.Range("AD5").Value = 11
$ws.Range("AB5").Value = 15
$ws.Range("AJ5").Value = 75
$ws.Range("V5").Value = 1.82
$ws.Range("AE5").Value = 23
$ws.Range("AL6").Value = 1000
$ws.Range("Z6").Value = 1000
$ws.Range("AN6").Value = 1000
$ws.Range("M6").Value = 1.11
$ws.Range("S6").Value = 5.3
$ws.Range("L6").Value = 1.01
$ws.Range("P6").Value = 1.37
$ws.Range("N6").Value = 1.37
$ws.Range("R6").Value = 1.11
$ws.Range("AK6").Value = 1000
$ws.Range("Q6").Value = 2.7
$ws.Range("AO6").Value = 1000
$ws.Range("AF6").Value = 1000
$ws.Range("X6").Value = 1000
$ws.Range("AC6").Value = 1000
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = 8.6
$ws.Range("G6").Value = 2.12
$ws.Range("Y6").Value = 1000
$ws.Range("O6").Value = 1.11
$ws.Range("AH6").Value = 1000
$ws.Range("AG6").Value = 1000
$ws.Range("K6").Value = 4.3
$ws.Range("AM6").Value = 1000
$ws.Range("J6").Value = 2.42
$ws.Range("AA6").Value = 1000
$ws.Range("U6").Value = 1.01
$ws.Range("F6").Value = 1.58
$ws.Range("V6").Value = 1.13
$ws.Range("AD6").Value = 1000
$ws.Range("AB6").Value = 1000
$ws.Range("T6").Value = 1.01
$ws.Range("AE6").Value = 1000
$ws.Range("AJ6").Value = 1000
$ws.Range("AI6").Value = 1000
$ws.Range("W6").Value = 1.89
$ws.Range("AO7").Value = 1000
$ws.Range("AK7").Value = 1000
$ws.Range("AJ7").Value = 1000
$ws.Range("L7").Value = 1.01
$ws.Range("P7").Value = 1.08
$ws.Range("U7").Value = 1.01
$ws.Range("AG7").Value = 1000
$ws.Range("AF7").Value = 1000
$ws.Range("O7").Value = 1.01
$ws.Range("Z7").Value = 1000
$ws.Range("M7").Value = 1.01
$ws.Range("X7").Value = 1000
$ws.Range("AD7").Value = 1000
$ws.Range("R7").Value = 1.08
$ws.Range("T7").Value = 1.01
$ws.Range("V7").Value = 1.01
$ws.Range("AC7").Value = 1000
$ws.Range("AI7").Value = 1000
$ws.Range("AN7").Value = 1000
$ws.Range("Y7").Value = 1000
$ws.Range("AM7").Value = 1000
$ws.Range("AH7").Value = 1000
$ws.Range("J7").Value = 1.02
$ws.Range("AB7").Value = 1000
$ws.Range("W7").Value = 1.01
$ws.Range("AA7").Value = 1000
$ws.Range("AE7").Value = 1000
$ws.Range("N7").Value = 1.08
$ws.Range("S7").Value = 1.01
$ws.Range("AL7").Value = 1000
$ws.Range("AF8").Value = 65
$ws.Range("F9").Value = 1.96
$ws.Range("Q9").Value = 1.76
$ws.Range("G10").Value = 3.05
$ws.Range("F10").Value = 2.88
$ws.Range("I10").Value = 2.64
$ws.Range("G12").Value = 1.88
$ws.Range("H12").Value = 4.4
$ws.Range("I12").Value = 4.8
$ws.Range("F12").Value = 1.86
$ws.Range("Q12").Value = 1.73
$ws.Range("F13").Value = 1.28
$ws.Range("G13").Value = 1.31
$ws.Range("J13").Value = 6
$ws.Range("O13").Value = 1.22
$ws.Range("P13").Value = 2.32
$ws.Range("S13").Value = 2.76
$ws.Range("T13").Value = 2.24
$ws.Range("Y13").Value = 95
$ws.Range("AD13").Value = 55
$ws.Range("G14").Value = 1.43
$ws.Range("J14").Value = 4.9
$ws.Range("H14").Value = 9
$ws.Range("F14").Value = 1.42
$ws.Range("K14").Value = 5.3
$ws.Range("P14").Value = 2.16
$ws.Range("F17").Value = 1.49
$ws.Range("R17").Value = 1.61
$ws.Range("AI17").Value = 80
$ws.Range("AH17").Value = 24
$ws.Range("AJ17").Value = 14.5
$ws.Range("AG17").Value = 10.5
$ws.Range("G18").Value = 1.93
$ws.Range("F18").Value = 1.92
$ws.Range("H18").Value = 4.1
$ws.Range("I19").Value = 6.6
$ws.Range("K19").Value = 4.5
$ws.Range("F19").Value = 1.63
$ws.Range("H19").Value = 5.9
$ws.Range("G19").Value = 1.67
$ws.Range("T19").Value = 1.93
$ws.Range("AB19").Value = 9
$ws.Range("Y19").Value = 21
$ws.Range("AD19").Value = 24
$ws.Range("Z19").Value = 140
$ws.Range("X19").Value = 17
$ws.Range("AJ19").Value = 16
$ws.Range("AH19").Value = 23
$ws.Range("K20").Value = 3.75
$ws.Range("J20").Value = 3.65
$ws.Range("J21").Value = 5.5
$ws.Range("K21").Value = 5.6
$ws.Range("I21").Value = 11
$ws.Range("O21").Value = 1.2
$ws.Range("U21").Value = 1.9
$ws.Range("T21").Value = 2
$ws.Range("Y21").Value = 70
$ws.Range("AD21").Value = 85
$ws.Range("AC21").Value = 12.5
$ws.Range("F22").Value = 1.68
$ws.Range("K22").Value = 4.5
$ws.Range("I22").Value = 5.7
$ws.Range("P22").Value = 2.5
$ws.Range("H23").Value = 4.8
$ws.Range("G24").Value = 2.06
$ws.Range("Q24").Value = 1.73
$ws.Range("I26").Value = 3.45
$ws.Range("H26").Value = 3.15
$ws.Range("F26").Value = 2.7
$ws.Range("G26").Value = 2.98
$ws.Range("F29").Value = 2.08
$ws.Range("H29").Value = 4.5
$ws.Range("G29").Value = 2.24
$ws.Range("Q29").Value = 2.8
$ws.Range("P29").Value = 1.45
$ws.Range("F30").Value = 2.14
$ws.Range("H30").Value = 4.4
$ws.Range("F31").Value = 1.84
$ws.Range("L31").Value = 1.33
$ws.Range("S31").Value = 3.4
$ws.Range("Q31").Value = 1.94
$ws.Range("P31").Value = 1.92
$ws.Range("AG31").Value = 10.5
